# Lab 12 part 2 -- re-knit of labs/illusion_key.docx
#
# Changes applied:
#   1. Drop the stray acute-accent character that used to trail
#      "...which itself is based on " right before the "Kohske Takahashi"
#      hyperlink.
#   2. Turn the separate " " + "code at" text that followed the
#      "Kohske Takahashi" hyperlink into "'s code at" (i.e. "Kohske
#      Takahashi's code at ...").
#   3. Bump the "Date this report was generated" timestamp.

$d = $word.ActiveDocument

# Make sure straight apostrophes stay straight (no curly-quote autocorrect)
# while we edit.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# --- 1) Remove the stray "´" after "which itself is based on " -----------
$rAccent = $d.Content
$rAccent.Find.Execute("based on ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$accentChar = $d.Range($rAccent.End, $rAccent.End + 1)
$accentChar.Delete()

# --- 2) "Kohske Takahashi code at" -> "Kohske Takahashi's code at" -------
# Insert "'s " right inside the plain ("code at") run so the new text does
# not pick up the hyperlink's character style, then drop the old space
# that used to sit between the hyperlink and "code".
$rCode = $d.Content
$rCode.Find.Execute("code at", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$cR = $d.Range($rCode.Start, $rCode.Start + 1)
$cR.InsertBefore("'s ")

$rHyper = $d.Content
$rHyper.Find.Execute("Kohske Takahashi", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$oldSpace = $d.Range($rHyper.End + 1, $rHyper.End + 2)
$oldSpace.Delete()

# --- 3) Update the generation timestamp -----------------------------------
$d.Content.Find.Execute( `
    "Date this report was generated: 2015-06-18 13:05:56.", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Date this report was generated: 2015-06-18 13:36:19.", 2)
